$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the original style of column D (price) so we can restore it after
# forcing a Text number format. Forcing Text avoids Excel re-interpreting the
# numeric-looking price strings (e.g. "316.35") as IEEE-754 doubles, which would
# otherwise be re-serialized with binary floating-point noise (e.g. "316.35000000000002").
$priceRange = $ws.Range("D2:D51")
$origStyle = $priceRange.Style
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "41.715.66"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "2.466.86"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "316.35"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "92.93"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "0.550"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  +3.54%  "
$ws.Range("D10").Value = "32.73"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  +7.48%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "2.848.64"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "6.90"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "15.78"
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("D16").Value = "2.463.12"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "0.782"
$ws.Range("E17").Value = "  +4.19%  "
$ws.Range("D18").Value = "41.705.63"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").Value = "11.53"
$ws.Range("E21").Value = "  +3.19%  "
$ws.Range("D22").Value = "71.11"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").Value = "239.81"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").Value = "2.73"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").Value = "1.93"
$ws.Range("E25").Value = "  +1.85%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "24.74"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("D29").Value = "9.80"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("D30").Value = "35.72"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("D31").Value = "156.07"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "5.51"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "2.51"
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("D36").Value = "17.55"
$ws.Range("E36").Value = "  -3.38%  "
$ws.Range("D37").Value = "2.88"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("D41").Value = "3.99"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "1.976.38"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").Value = "19.01"
$ws.Range("E44").Value = "  -4.74%  "
$ws.Range("D45").Value = "0.0285"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").Value = "9.05"
$ws.Range("E47").Value = "  +2.91%  "
$ws.Range("D48").Value = "2.701.79"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").Value = "67.02"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").Value = "72.88"
$ws.Range("E51").Value = "  -0.57%  "

# Restore the original (unformatted) style now that the text values are set.
$priceRange.Style = $origStyle
